$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain numeric-looking price text to stay text (matches source data,
# which stores these as literal strings, not localized numbers).
foreach ($addr in @("D5", "D6", "D10", "D14", "D17", "D18", "D20", "D21", "D23", "D24", "D25", "D27", "D28", "D30", "D32", "D33", "D34", "D36", "D40", "D41", "D43", "D45", "D46", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.544.07"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "3.105.27"

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "386.28"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("D6").Value = "103.31"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D10").Value = "37.10"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").Value = "3.592.97"
$ws.Range("E13").Value = "  +2.90%  "

$ws.Range("D14").Value = "18.58"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D16").Value = "3.097.04"
$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("D17").Value = "0.996"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").Value = "10.99"
$ws.Range("E18").Value = "  +5.37%  "

$ws.Range("D19").Value = "51.574.94"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").Value = "3.27"
$ws.Range("E20").Value = "  +7.84%  "

$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").Value = "69.91"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "266.79"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("E26").Value = "  -1.58%  "

$ws.Range("D27").Value = "27.02"
$ws.Range("E27").Value = "  +3.30%  "

$ws.Range("D28").Value = "7.25"
$ws.Range("E28").Value = "  -3.63%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").Value = "0.166"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("D32").Value = "10.36"
$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").Value = "0.0482"
$ws.Range("E33").Value = "  +5.55%  "

$ws.Range("D34").Value = "35.27"
$ws.Range("E34").Value = "  +3.46%  "

$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").Value = "50.07"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("E38").Value = "  +2.03%  "

$ws.Range("E39").Value = "  +2.75%  "

$ws.Range("D40").Value = "1.88"
$ws.Range("E40").Value = "  +1.53%  "

$ws.Range("D41").Value = "129.01"
$ws.Range("E41").Value = "  +1.86%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").Value = "16.57"
$ws.Range("E43").Value = "  -3.72%  "

$ws.Range("E44").Value = "  -3.65%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.70"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "22.43"
$ws.Range("E46").Value = "  +4.02%  "

$ws.Range("E47").Value = "  +4.72%  "

$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").Value = "2.068.27"
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  +3.96%  "

$ws.Range("D51").Value = "0.932"
$ws.Range("E51").Value = "  +17.39%  "
